$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.985.77'
$ws.Range('E2').Value = '  +4.42%  '
$ws.Range('D3').Value = '2.779.34'
$ws.Range('E3').Value = '  +4.81%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.64%  '
$ws.Range('E6').Value = '  +2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.63%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0857'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.99'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = '3.218.44'
$ws.Range('E15').Value = '  +4.97%  '
$ws.Range('D16').Value = '2.783.20'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('D17').Value = '51.894.97'
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.878'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '0.0₃0980'
$ws.Range('E22').Value = '  +2.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '276.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0819'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0386'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.56%  '
$ws.Range('E41').Value = '  +25.90%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.17%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.116'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '23.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '2.070.64'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.895'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.73%  '
